# Updated symbol list pricing/volume data (per commit: "Updated symbol list on Tue Jan 24 21:56:19 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D = Price, E = Volume(1h) -- both columns hold text-formatted numeric/percent
# strings (e.g. "35.50", "0.40%"), so force the Text number format before
# assigning the value. Otherwise Excel auto-converts the numeric-looking text
# into a real number/percentage and silently drops the exact text (e.g. "35.50"
# -> 35.5, "2.900" -> 2.9, "0.40%" -> 0.004).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $ws.Range("D2") "307.37"
Set-TextValue $ws.Range("E2") "0.40%"
Set-TextValue $ws.Range("D3") "35.50"
Set-TextValue $ws.Range("E3") "-2.21%"
Set-TextValue $ws.Range("E4") "1.09%"
Set-TextValue $ws.Range("D5") "0.08152"
Set-TextValue $ws.Range("E5") "2.76%"
Set-TextValue $ws.Range("D6") "1.970"
Set-TextValue $ws.Range("E6") "-11.97%"
Set-TextValue $ws.Range("D7") "7.946"
Set-TextValue $ws.Range("E7") "-0.76%"
Set-TextValue $ws.Range("D8") "2.897"
Set-TextValue $ws.Range("E8") "10.08%"
Set-TextValue $ws.Range("D9") "0.9280"
Set-TextValue $ws.Range("E9") "-0.01%"
Set-TextValue $ws.Range("D10") "0.1084"
Set-TextValue $ws.Range("E10") "10.45%"
Set-TextValue $ws.Range("D11") "0.1928"
Set-TextValue $ws.Range("E11") "2.55%"
Set-TextValue $ws.Range("D12") "0.09644"
Set-TextValue $ws.Range("E12") "6.32%"
Set-TextValue $ws.Range("D13") "0.03634"
Set-TextValue $ws.Range("E13") "-2.15%"
Set-TextValue $ws.Range("D14") "0.09915"
Set-TextValue $ws.Range("E14") "-0.09%"
Set-TextValue $ws.Range("D15") "0.001437"
Set-TextValue $ws.Range("E15") "0.15%"
Set-TextValue $ws.Range("D16") "0.005802"
Set-TextValue $ws.Range("E16") "3.51%"
Set-TextValue $ws.Range("D17") "3.478"
Set-TextValue $ws.Range("E17") "0.42%"
Set-TextValue $ws.Range("D18") "4.134"
Set-TextValue $ws.Range("E18") "-0.32%"
Set-TextValue $ws.Range("D19") "0.3418"
Set-TextValue $ws.Range("E19") "1.39%"
Set-TextValue $ws.Range("D20") "0.1312"
Set-TextValue $ws.Range("E20") "-0.51%"
Set-TextValue $ws.Range("D21") "5.131"
Set-TextValue $ws.Range("E21") "0.70%"
Set-TextValue $ws.Range("E22") "-2.51%"
Set-TextValue $ws.Range("D23") "0.04549"
Set-TextValue $ws.Range("E23") "0.14%"
Set-TextValue $ws.Range("D24") "0.001230"
Set-TextValue $ws.Range("E24") "-0.79%"
Set-TextValue $ws.Range("D25") "0.004764"
Set-TextValue $ws.Range("E25") "-0.28%"
Set-TextValue $ws.Range("D26") "0.0001259"
Set-TextValue $ws.Range("E26") "-3.09%"
Set-TextValue $ws.Range("D27") "0.0004458"
Set-TextValue $ws.Range("E27") "-5.91%"
Set-TextValue $ws.Range("D39") "0.01941"
Set-TextValue $ws.Range("E39") "1.26%"
Set-TextValue $ws.Range("D40") "0.04881"
Set-TextValue $ws.Range("E40") "-1.13%"
Set-TextValue $ws.Range("D41") "0.007839"
Set-TextValue $ws.Range("D42") "0.009697"
Set-TextValue $ws.Range("E42") "24.31%"
Set-TextValue $ws.Range("D43") "0.1380"
Set-TextValue $ws.Range("E43") "-1.19%"
Set-TextValue $ws.Range("D44") "0.002127"
Set-TextValue $ws.Range("E44") "0.10%"
Set-TextValue $ws.Range("D45") "0.01156"
Set-TextValue $ws.Range("D46") "0.00006510"
Set-TextValue $ws.Range("E46") "4.45%"
Set-TextValue $ws.Range("D47") "0.00000000752"
Set-TextValue $ws.Range("E47") "0.17%"
Set-TextValue $ws.Range("E48") "23.12%"
Set-TextValue $ws.Range("D49") "0.001303"
Set-TextValue $ws.Range("E49") "-27.65%"
Set-TextValue $ws.Range("D50") "0.00002105"
Set-TextValue $ws.Range("E50") "0.17%"
Set-TextValue $ws.Range("D51") "0.0002005"
Set-TextValue $ws.Range("E51") "0.17%"
